$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row correct-answer marks value (B11: 3 -> 5)
$ws.Range("B11").Value = 5

# Update "Total" row values (B12: 84 -> 140) and the corresponding
# "Corr/total" display text in E12 (84/84 -> 140/140)
$ws.Range("B12").Value = 140
$ws.Range("E12").Value = "140/140"
